$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 1 de Agosto de 2020 a las 15:32'

# Row 4
$ws.Cells.Item(4, 2).Value = 4707584
$ws.Cells.Item(4, 3).Value = 1695
$ws.Cells.Item(4, 5).Value = 2222361

# Row 6
$ws.Cells.Item(6, 2).Value = 1722159
$ws.Cells.Item(6, 3).Value = 25105
$ws.Cells.Item(6, 4).Value = 1121105
$ws.Cells.Item(6, 5).Value = 564224
$ws.Cells.Item(6, 7).Value = 279
$ws.Cells.Item(6, 8).Value = 36830

# Row 14
$ws.Cells.Item(14, 2).Value = 303952
$ws.Cells.Item(14, 3).Value = 771

# Row 17
$ws.Cells.Item(17, 2).Value = 277478
$ws.Cells.Item(17, 3).Value = 1573
$ws.Cells.Item(17, 4).Value = 237548
$ws.Cells.Item(17, 5).Value = 37043
$ws.Cells.Item(17, 7).Value = 21
$ws.Cells.Item(17, 8).Value = 2887

# Row 41
$ws.Cells.Item(41, 2).Value = 67448
$ws.Cells.Item(41, 3).Value = 491
$ws.Cells.Item(41, 4).Value = 58525
$ws.Cells.Item(41, 5).Value = 8470
$ws.Cells.Item(41, 7).Value = 6
$ws.Cells.Item(41, 8).Value = 453

# Row 44
$ws.Cells.Item(44, 2).Value = 54732
$ws.Cells.Item(44, 3).Value = 431
$ws.Cells.Item(44, 7).Value = 1
$ws.Cells.Item(44, 8).Value = 6148

# Row 59
$ws.Cells.Item(59, 2).Value = 32157
$ws.Cells.Item(59, 3).Value = 279
$ws.Cells.Item(59, 4).Value = 26474
$ws.Cells.Item(59, 5).Value = 5229
$ws.Cells.Item(59, 7).Value = 6
$ws.Cells.Item(59, 8).Value = 454

# Row 62
$ws.Cells.Item(62, 2).Value = 25882
$ws.Cells.Item(62, 3).Value = 330
$ws.Cells.Item(62, 5).Value = 11253
$ws.Cells.Item(62, 7).Value = 9
$ws.Cells.Item(62, 8).Value = 582

# Row 64
$ws.Cells.Item(64, 1).Value = 'Uzbekistan'
$ws.Cells.Item(64, 2).Value = 24569
$ws.Cells.Item(64, 3).Value = 560
$ws.Cells.Item(64, 4).Value = 14916
$ws.Cells.Item(64, 5).Value = 9508
$ws.Cells.Item(64, 7).Value = 4
$ws.Cells.Item(64, 8).Value = 145

# Row 65
$ws.Cells.Item(65, 1).Value = 'Marruecos'
$ws.Cells.Item(65, 2).Value = 24322
$ws.Cells.Item(65, 3).Value = 0
$ws.Cells.Item(65, 4).Value = 17658
$ws.Cells.Item(65, 5).Value = 6311
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 8).Value = 353

# Row 79
$ws.Cells.Item(79, 5).Value = 6753
$ws.Cells.Item(79, 7).Value = 1
$ws.Cells.Item(79, 8).Value = 83

# Row 143
$ws.Cells.Item(143, 1).Value = 'Uganda'
$ws.Cells.Item(143, 2).Value = 1176
$ws.Cells.Item(143, 3).Value = 22
$ws.Cells.Item(143, 4).Value = 1045
$ws.Cells.Item(143, 5).Value = 127
$ws.Cells.Item(143, 7).Value = 1
$ws.Cells.Item(143, 8).Value = 4

# Row 144
$ws.Cells.Item(144, 1).Value = 'Georgia'
$ws.Cells.Item(144, 2).Value = 1171
$ws.Cells.Item(144, 3).Value = 3
$ws.Cells.Item(144, 4).Value = 947
$ws.Cells.Item(144, 5).Value = 207
$ws.Cells.Item(144, 8).Value = 17

# Row 179
$ws.Cells.Item(179, 2).Value = 188
$ws.Cells.Item(179, 3).Value = 1
$ws.Cells.Item(179, 4).Value = 182

# Row 190
$ws.Cells.Item(190, 1).Value = 'Papua Nueva Guinea'
$ws.Cells.Item(190, 3).Value = 19
$ws.Cells.Item(190, 4).Value = 34
$ws.Cells.Item(190, 5).Value = 55
$ws.Cells.Item(190, 8).Value = 2

# Row 191
$ws.Cells.Item(191, 1).Value = 'Antigua y Barbuda'
$ws.Cells.Item(191, 2).Value = 91
$ws.Cells.Item(191, 4).Value = 67
$ws.Cells.Item(191, 5).Value = 21
$ws.Cells.Item(191, 8).Value = 3

# Row 192
$ws.Cells.Item(192, 1).Value = 'Liechtenstein'
$ws.Cells.Item(192, 2).Value = 88
$ws.Cells.Item(192, 4).Value = 85
$ws.Cells.Item(192, 5).Value = 2
$ws.Cells.Item(192, 8).Value = 1
